$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 4.2
$ws.Range("I2").Value = 1.75
$ws.Range("J2").Value = 4.33
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 29
$ws.Range("AK2").Value = 15
$ws.Range("AN2").Value = 6.5
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 23
$ws.Range("AS2").Value = 126
$ws.Range("AX2").Value = 4

# Row 3 updates
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 2.9
$ws.Range("J3").Value = 3.5
$ws.Range("AK3").Value = 29
$ws.Range("AN3").Value = 4.5
$ws.Range("AO3").Value = 17

# Row 4 updates
$ws.Range("G4").Value = 1.85
$ws.Range("I4").Value = 4.1
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("Z4").Value = 15
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 15
$ws.Range("AU4").Value = 7.5

# Row 5 updates
$ws.Range("G5").Value = 2.88
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 2.25
$ws.Range("J5").Value = 3.4
$ws.Range("U5").Value = 1.55
$ws.Range("V5").Value = 2.15
$ws.Range("X5").Value = 15
$ws.Range("Z5").Value = 29
$ws.Range("AJ5").Value = 9.5
$ws.Range("AO5").Value = 15
$ws.Range("AP5").Value = 21
$ws.Range("AZ5").Value = 21

# Row 6 updates
$ws.Range("G6").Value = 1.57
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 2.2
$ws.Range("K6").Value = 2.2
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.75
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("AN6").Value = 3.4
$ws.Range("AO6").Value = 8
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 26
$ws.Range("AU6").Value = 9.5
$ws.Range("AX6").Value = 7
$ws.Range("AY6").Value = 34
$ws.Range("BA6").Value = 126

# Row 7 updates
$ws.Range("G7").Value = 5.25
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 1.6
$ws.Range("J7").Value = 5
$ws.Range("L7").Value = 2.1
$ws.Range("U7").Value = 1.57
$ws.Range("V7").Value = 2.25
$ws.Range("W7").Value = 19
$ws.Range("X7").Value = 29
$ws.Range("Y7").Value = 17
$ws.Range("Z7").Value = 51
$ws.Range("AB7").Value = 34
$ws.Range("AD7").Value = 8
$ws.Range("AE7").Value = 13
$ws.Range("AG7").Value = 126
$ws.Range("AH7").Value = 10
$ws.Range("AK7").Value = 13
$ws.Range("AN7").Value = 7
$ws.Range("AO7").Value = 23
$ws.Range("AR7").Value = 81
$ws.Range("AU7").Value = 7.5
$ws.Range("AX7").Value = 4
$ws.Range("AY7").Value = 8
$ws.Range("BA7").Value = 23
$ws.Range("BC7").Value = 81

# Row 8 updates
$ws.Range("G8").Value = 1.13
$ws.Range("H8").Value = 8.5
$ws.Range("I8").Value = 19
$ws.Range("K8").Value = 3.25
$ws.Range("L8").Value = 12
$ws.Range("O8").Value = 1.08
$ws.Range("P8").Value = 8
$ws.Range("Q8").Value = 1.3
$ws.Range("R8").Value = 3.5
$ws.Range("S8").Value = 1.18
$ws.Range("T8").Value = 4.5
$ws.Range("AC8").Value = 23
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 29
$ws.Range("AG8").Value = 351
$ws.Range("AM8").Value = 67
$ws.Range("AQ8").Value = 9.5
$ws.Range("AT8").Value = 4.5
$ws.Range("AX8").Value = 15
$ws.Range("BA8").Value = 351

# Row 14 updates
$ws.Range("G14").Value = 28
$ws.Range("H14").Value = 9.75
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 4.35
$ws.Range("L14").Value = 1.2
$ws.Range("O14").Value = 1.04
$ws.Range("P14").Value = 8.25
$ws.Range("Q14").Value = 1.16
$ws.Range("R14").Value = 4.55
$ws.Range("S14").Value = 1.1
$ws.Range("T14").Value = 5.7
$ws.Range("U14").Value = 2.37
$ws.Range("V14").Value = 1.52
$ws.Range("W14").Value = 175
$ws.Range("X14").Value = 800
$ws.Range("Y14").Value = 175
$ws.Range("AB14").Value = 500
$ws.Range("AC14").Value = 35
$ws.Range("AD14").Value = 35
$ws.Range("AE14").Value = 70
$ws.Range("AF14").Value = 300
$ws.Range("AH14").Value = 17.5
$ws.Range("AI14").Value = 9.25
$ws.Range("AJ14").Value = 18
$ws.Range("AK14").Value = 7
$ws.Range("AL14").Value = 13.5
$ws.Range("AN14").Value = 32
$ws.Range("AO14").Value = 250
$ws.Range("AP14").Value = 120
$ws.Range("AT14").Value = 5.7
$ws.Range("AU14").Value = 14.5
$ws.Range("AX14").Value = 3.8
$ws.Range("AY14").Value = 3.9
$ws.Range("BA14").Value = 5.7
$ws.Range("BB14").Value = 21
$ws.Range("BC14").Value = 150

# Row 15 updates
$ws.Range("M15").Value = 1.05
$ws.Range("N15").Value = 11

# Row 17 updates
$ws.Range("AX17").Value = 5
